$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '88.099.32'
$ws.Range('E2').Value = '  +10.82%  '
$ws.Range('D3').Value = '3.341.79'
$ws.Range('E3').Value = '  +6.40%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '219.50'
$ws.Range('E5').Value = '  +6.68%  '
$ws.Range('D6').Value = '638.96'
$ws.Range('E6').Value = '  +2.79%  '
$ws.Range('D7').Value = '0.328'
$ws.Range('E7').Value = '  +24.50%  '
$ws.Range('E8').Value = '  -0.11%  '
$ws.Range('E9').Value = '  +6.47%  '
$ws.Range('D10').Value = '3.350.46'
$ws.Range('E10').Value = '  +6.77%  '
$ws.Range('D11').Value = '0.613'
$ws.Range('E11').Value = '  +5.90%  '
$ws.Range('D12').Value = '0.0000278'
$ws.Range('E12').Value = '  +12.07%  '
$ws.Range('E13').Value = '  +2.33%  '
$ws.Range('D14').Value = '3.970.05'
$ws.Range('E14').Value = '  +6.78%  '
$ws.Range('D15').Value = '34.54'
$ws.Range('E15').Value = '  +10.65%  '
$ws.Range('D16').Value = '5.43'
$ws.Range('E16').Value = '  +4.52%  '
$ws.Range('D17').Value = '87.831.65'
$ws.Range('E17').Value = '  +10.66%  '
$ws.Range('D18').Value = '3.354.80'
$ws.Range('E18').Value = '  +7.07%  '
$ws.Range('D19').Value = '14.66'
$ws.Range('E19').Value = '  +3.97%  '
$ws.Range('D20').Value = '3.21'
$ws.Range('E20').Value = '  +9.47%  '
$ws.Range('D21').Value = '450.05'
$ws.Range('D22').Value = '9.14'
$ws.Range('E22').Value = '  +1.62%  '
$ws.Range('D23').Value = '5.34'
$ws.Range('E23').Value = '  +4.06%  '
$ws.Range('E24').Value = '  +7.06%  '
$ws.Range('D25').Value = '5.41'
$ws.Range('E25').Value = '  +17.14%  '
$ws.Range('D26').Value = '12.34'
$ws.Range('E26').Value = '  +15.58%  '
$ws.Range('D27').Value = '3.520.59'
$ws.Range('E27').Value = '  +6.72%  '
$ws.Range('D28').Value = '78.94'
$ws.Range('E28').Value = '  +4.86%  '
$ws.Range('E29').Value = '  +9.25%  '
$ws.Range('E30').Value = '  +0.05%  '
$ws.Range('E31').Value = '  +56.29%  '
$ws.Range('D32').Value = '604.96'
$ws.Range('E32').Value = '  +10.26%  '
$ws.Range('D33').Value = '9.33'
$ws.Range('E33').Value = '  +5.56%  '
$ws.Range('D34').Value = '0.998'
$ws.Range('E34').Value = '  -0.16%  '
$ws.Range('E35').Value = '  +7.13%  '
$ws.Range('E36').Value = '  +4.48%  '
$ws.Range('D37').Value = '0.152'
$ws.Range('E37').Value = '  +2.85%  '
$ws.Range('D38').Value = '23.49'
$ws.Range('E38').Value = '  +3.40%  '
$ws.Range('D39').Value = '6.76'
$ws.Range('E39').Value = '  +22.86%  '
$ws.Range('D40').Value = '0.421'
$ws.Range('E40').Value = '  +6.09%  '
$ws.Range('D41').Value = '0.999'
$ws.Range('E41').Value = '  -0.04%  '
$ws.Range('E42').Value = '  +3.13%  '
$ws.Range('E43').Value = '  +16.20%  '
$ws.Range('D44').Value = '3.09'
$ws.Range('E44').Value = '  +16.40%  '
$ws.Range('D45').Value = '158.25'
$ws.Range('E45').Value = '  -2.90%  '
$ws.Range('E46').Value = '  +0.03%  '
$ws.Range('D47').Value = '190.41'
$ws.Range('E47').Value = '  +1.57%  '
$ws.Range('B48').Value = 'OKB'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D48').Value = '46.16'
$ws.Range('E48').Value = '  +8.61%  '
$ws.Range('B49').Value = 'ImmutableX'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D49').Value = '1.38'
$ws.Range('E49').Value = '  +8.69%  '
$ws.Range('D50').Value = '0.789'
$ws.Range('E50').Value = '  +1.38%  '
$ws.Range('D51').Value = '26.75'
$ws.Range('E51').Value = '  +10.43%  '
